# Add a new row (row 12) of data for year 2021 to Sheet1, following the
# same pattern/formatting as the existing data rows (2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A: year label, copy an existing "year" cell first so the new cell
# picks up the exact same style (border/bold/centered) used by A2:A11,
# then overwrite its text.
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "2021年"

# Numeric columns with values for 2021.
$ws.Range("B12").Value = 0.03
$ws.Range("E12").Value = 13.7
$ws.Range("F12").Value = 1.42
$ws.Range("G12").Value = 0.36
$ws.Range("H12").Value = 9.18
$ws.Range("I12").Value = 11.32
$ws.Range("J12").Value = 2.98
$ws.Range("L12").Value = 13.7
$ws.Range("N12").Value = 13.34
$ws.Range("S12").Value = 1.51
$ws.Range("V12").Value = 0.96

# Columns with no data for 2021 (C, D, K, M, O, P, Q, R, T, U). Copy one of
# the existing empty cells (e.g. C11) into each of these so that an empty
# cell entry is still created at that position, matching the layout of the
# other rows instead of leaving the cell completely absent.
$ws.Range("C11").Copy($ws.Range("C12"))
$ws.Range("C11").Copy($ws.Range("D12"))
$ws.Range("C11").Copy($ws.Range("K12"))
$ws.Range("C11").Copy($ws.Range("M12"))
$ws.Range("C11").Copy($ws.Range("O12"))
$ws.Range("C11").Copy($ws.Range("P12"))
$ws.Range("C11").Copy($ws.Range("Q12"))
$ws.Range("C11").Copy($ws.Range("R12"))
$ws.Range("C11").Copy($ws.Range("T12"))
$ws.Range("C11").Copy($ws.Range("U12"))
